# Tracking workbook update:
#  - add new tracking entry "dokoncit sprint burdowncharts (sprint1-6)" in row 6 (col B)
#  - highlight the updated task cell (B2) and the new cell (B6) with the
#    accent6 theme color fill
#  - give B6 a left/right thin border (matching the "new entry" look)
#  - widen column B slightly to fit the new text
#  - leave the selection on B3 (matches the author's final selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) New row with the new sprint task description
$ws.Range("B6").Value = "dokoncit sprint burdowncharts (sprint1-6)"

# 2) Style the new B6 cell: accent6 theme fill, with only left/right thin
#    borders (no top/bottom) to set it apart as a freshly appended row.
$b6 = $ws.Range("B6")
$b6.Interior.ThemeColor = 10
$b6.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeLeft).LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous
$b6.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeLeft).Weight = [Microsoft.Office.Interop.Excel.XlBorderWeight]::xlThin
$b6.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeRight).LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous
$b6.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeRight).Weight = [Microsoft.Office.Interop.Excel.XlBorderWeight]::xlThin

# 3) Highlight the edited "Use case description update" cell (B2) with the
#    same accent6 theme color fill, keeping its existing full thin border.
$b2 = $ws.Range("B2")
$b2.Interior.ThemeColor = 10

# 4) Widen column B a bit so the new (longer) text fits nicely.
$ws.Columns.Item(2).ColumnWidth = 38

# 5) Final selection left on B3, as in the saved workbook.
$ws.Range("B3").Select()

# 6) Restore the (smaller/maximized) window size & position used on the
#    author's machine when the file was last saved.
$win = $wb.Windows.Item(1)
$win.Left = 3465
$win.Top = 3465
$win.Width = 21600
$win.Height = 11385
